# Adds a new "Creamos un nuevo servicio" / "ng g s {nombre_servicio}" section
# to the end of the cheat-sheet (mirrors the existing "Creamos nueva
# directiva" section), and removes a stray leftover test value (B36 = 44).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the stray leftover numeric test value in B36.
$ws.Cells.Item(36, 2).ClearContents()

# Insert a new border/separator row at 35, reusing the same formatting as
# the existing separator row 29 (columns C:G only, so column B stays
# unstyled, matching the pattern used by every other separator row).
$ws.Range("C29:G29").Copy()
$ws.Range("C35:G35").PasteSpecial(-4122)

# New bold section header in row 37.
$ws.Cells.Item(37, 2).Value = "Creamos un nuevo servicio"
$ws.Cells.Item(37, 2).Font.Bold = $true

# New content/example row in row 39.
$ws.Cells.Item(39, 2).Value = "ng g s {nombre_servicio}"

# Update the sheet view to reflect the new scroll position / selection.
$null = $ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$null = $ws.Range("B42").Select()
